# Add files via upload
# Populates the "standard code" (표준코드) columns on Sheet1 with newly
# digitised values, and appends three new drug rows (B1CAG11 / B1AMD11 /
# B1INDC31) that were missing from the master list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fill in previously-blank 표준코드 (standard code) cells -------------
$ws.Range("C8").Value  = "06449036011"   # Heparin 20000iu 20mL(중외)
$ws.Range("C10").Value = "06789008617"   # K20mEq/NS 100mL(NaK주 200)중외
$ws.Range("C11").Value = "06451008631"   # KCl-40 20mL(대한)
$ws.Range("C13").Value = "06963004015"   # Leuplin DPS 11.25mg(다케다)
$ws.Range("C17").Value = "06451011013"   # Magnesium Sulfate주 50% 20mL(대한)
$ws.Range("C21").Value = "06706075715"   # Pine(Heparin) 5000iu/5mL(휴온스)

# --- Append new drug rows --------------------------------------------------
$ws.Range("A26").Value = "B1CAG11"
$ws.Range("B26").Value = "Cal gluconate 10% 20mL(중외)"
$ws.Range("C26").Value = "06449142811"

$ws.Range("A27").Value = "B1AMD11"
$ws.Range("B27").Value = "CorDARONE 150mg(사노피)"
$ws.Range("C27").Value = "06520004519"

$ws.Range("A28").Value = "B1INDC31"
$ws.Range("B28").Value = "Carmine 0.8% 5mL(유나이티드)"
$ws.Range("C28").Value = "06443038516"

# --- Restore the cursor position left by the editor -------------------------
$ws.Range("D38").Select()
